$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.285.78"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "3.543.66"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'597.29"
$ws.Range("E5").Value = "  +2.61%  "
$ws.Range("D6").Value = "'138.15"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").Value = "3.541.82"
$ws.Range("E7").Value = "  +3.80%  "
$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  +3.55%  "
$ws.Range("D11").Value = "'6.91"
$ws.Range("E11").Value = "  -3.86%  "
$ws.Range("D12").Value = "'0.387"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("D13").Value = "4.141.90"
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("D14").Value = "'0.0000184"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "'27.32"
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("D16").Value = "3.542.12"
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "65.192.92"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  +4.69%  "
$ws.Range("D20").Value = "'5.92"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").Value = "'14.33"
$ws.Range("E21").Value = "  +5.00%  "
$ws.Range("D22").Value = "'392.55"
$ws.Range("E22").Value = "  +2.56%  "
$ws.Range("D23").Value = "'0.574"
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("D24").Value = "3.680.62"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("D25").Value = "'73.58"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'0.0000114"
$ws.Range("E27").Value = "  +8.44%  "
$ws.Range("D28").Value = "'7.89"
$ws.Range("E28").Value = "  +11.42%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'2.28"
$ws.Range("E30").Value = "  +3.48%  "
$ws.Range("D31").Value = "'8.24"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "3.562.38"
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'23.80"
$ws.Range("E34").Value = "  +4.85%  "
$ws.Range("D35").Value = "'0.145"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'1.33"
$ws.Range("E36").Value = "  +15.54%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.59"
$ws.Range("E37").Value = "  +8.82%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'170.43"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "'6.95"
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("D40").Value = "'5.02"
$ws.Range("E40").Value = "  +6.60%  "
$ws.Range("D41").Value = "'0.0806"
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("D42").Value = "'0.825"
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("D43").Value = "'26.52"
$ws.Range("E43").Value = "  +19.73%  "
$ws.Range("D44").Value = "'42.64"
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'4.44"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").Value = "'1.21"
$ws.Range("E47").Value = "  +9.75%  "
$ws.Range("D48").Value = "'1.68"
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("D49").Value = "'6.82"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("D50").Value = "2.381.01"
$ws.Range("E50").Value = "  +9.33%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'306.99"
$ws.Range("E51").Value = "  +10.88%  "
